$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New B:E (and derived G = sum(B:E)) values per row, regenerated s_vals
# data to filter save games. F (Win) column is untouched.
$data = @{
    2  = @(1.459612070389937,  1.667794583268128,  0.1575252929769615, 0.496779210170732,  3.781711156805759)
    3  = @(1.459612070389937,  1.667794583268128,  0.8054896365839992, 8.660232485948974,  12.59312877619104)
    4  = @(3.230985683306322,  1.667794583268128,  26.21740644021617,  0.496779210170732,  31.61296591696135)
    5  = @(1.459612070389937,  1.667794583268128,  0.8054896365839992, 0.496779210170732,  4.429675500412797)
    6  = @(3.230985683306322,  1.667794583268128,  0.1575252929769615, 0.496779210170732,  5.553084769722144)
    7  = @(0.01514828764759746,0.002777888934908601,0.1575252929769615,0.496779210170732,  0.6722306797301996)
    8  = @(3.230985683306322,  1.667794583268128,  0.8054896365839992, 0.496779210170732,  6.201049113329182)
    9  = @(3.230985683306322,  1.667794583268128,  0.1575252929769615, 0.496779210170732,  5.553084769722144)
    10 = @(3.230985683306322,  1.667794583268128,  0.1575252929769615, 0.496779210170732,  5.553084769722144)
    11 = @(3.230985683306322,  1.667794583268128,  0.1575252929769615, 8.660232485948974,  13.71653804550039)
    12 = @(3.230985683306322,  1.667794583268128,  0.8054896365839992, 0.496779210170732,  6.201049113329182)
    13 = @(0.0001488876196638067,0.002777888934908601,0.1575252929769615,0.496779210170732,0.6572312797022659)
    14 = @(3.230985683306322,  1.667794583268128,  0.8054896365839992, 0.496779210170732,  6.201049113329182)
    15 = @(0.3048080303191223, 0.04240448674262143, 0.8054896365839992, 0.496779210170732,  1.649481363816475)
    16 = @(0.3048080303191223, 0.3127903958511391,  0.1575252929769615, 0.496779210170732,  1.271902929317955)
    17 = @(0.127881588408715,  0.3127903958511391,  3.900430680208489,  8.660232485948974,  13.00133515041732)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 7).Value = $vals[4]
}
